# Auto-generated Excel COM-interop script implementing the crime-data update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report date range) ---
$volCell = $ws.Range("A8")
$volCell.Characters(21, 2).Text = "35"

$dateCell = $ws.Range("C9")
$dateCell.Characters(27, 9).Text = "8/26/2024"
$dateCell.Characters(47, 9).Text = "9/1/2024"

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("F14").Value = 2
$ws.Range("I14").Value = 2
$ws.Range("J14").Value = 3
$ws.Range("K14").Value = -33.333333333333
$ws.Range("M15").Value = -57.142857142857
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 83
$ws.Range("J16").Value = 74
$ws.Range("K16").Value = 12.162162162162
$ws.Range("L16").Value = 2.469135802469
$ws.Range("M16").Value = -44.666666666666
$ws.Range("N16").Value = -81.758241758241
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 21
$ws.Range("H17").Value = 75
$ws.Range("I17").Value = 117
$ws.Range("J17").Value = 111
$ws.Range("K17").Value = 5.405405405405
$ws.Range("L17").Value = 1.739130434782
$ws.Range("M17").Value = 42.682926829268
$ws.Range("N17").Value = -47.297297297297
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -66.666666666666
$ws.Range("I18").Value = 43
$ws.Range("J18").Value = 66
$ws.Range("K18").Value = -34.848484848484
$ws.Range("L18").Value = -44.155844155844
$ws.Range("M18").Value = -76.373626373626
$ws.Range("N18").Value = -94.49423815621
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 56
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = 1.818181818181
$ws.Range("I19").Value = 390
$ws.Range("J19").Value = 442
$ws.Range("K19").Value = -11.764705882352
$ws.Range("L19").Value = 16.071428571428
$ws.Range("M19").Value = 24.203821656051
$ws.Range("N19").Value = -12.751677852349
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 77.777777777777
$ws.Range("I20").Value = 120
$ws.Range("J20").Value = 84
$ws.Range("K20").Value = 42.857142857142
$ws.Range("L20").Value = 48.148148148148
$ws.Range("M20").Value = 5.263157894736
$ws.Range("N20").Value = -93.610223642172
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 30
$ws.Range("F21").Value = 110
$ws.Range("G21").Value = 98
$ws.Range("H21").Value = 12.244897959183
$ws.Range("I21").Value = 761
$ws.Range("J21").Value = 788
$ws.Range("K21").Value = -3.426395939086
$ws.Range("L21").Value = 9.182209469153
$ws.Range("M21").Value = -11.201866977829
$ws.Range("N21").Value = -80.041961709939
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -50
$ws.Range("J23").Value = 18
$ws.Range("K23").Value = 22.222222222222
$ws.Range("L23").Value = -15.384615384615
$ws.Range("M23").Value = 0
$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 85.714285714285
$ws.Range("F24").Value = 134
$ws.Range("G24").Value = 107
$ws.Range("H24").Value = 25.233644859813
$ws.Range("I24").Value = 1071
$ws.Range("J24").Value = 834
$ws.Range("K24").Value = 28.41726618705
$ws.Range("L24").Value = 42.231075697211
$ws.Range("M24").Value = 60.81081081081
$ws.Range("C25").Value = 34
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = 126.666666666667
$ws.Range("F25").Value = 112
$ws.Range("G25").Value = 86
$ws.Range("H25").Value = 30.232558139534
$ws.Range("I25").Value = 888
$ws.Range("J25").Value = 593
$ws.Range("K25").Value = 49.747048903878
$ws.Range("L25").Value = 84.232365145228
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = -55.555555555555
$ws.Range("F26").Value = 25
$ws.Range("H26").Value = 4.166666666666
$ws.Range("I26").Value = 216
$ws.Range("J26").Value = 185
$ws.Range("K26").Value = 16.756756756756
$ws.Range("L26").Value = 43.046357615894
$ws.Range("M26").Value = -11.475409836065
$ws.Range("C28").Value = 1
$ws.Range("I28").Value = 26
$ws.Range("K28").Value = 30
$ws.Range("L28").Value = 18.181818181818
$ws.Range("L29").Value = -81.818181818181
$ws.Range("L30").Value = -71.428571428571

# --- Numeric updates that also require a style/number-format correction ---
# (use Copy + PasteSpecial(xlPasteFormats) from a donor cell that already has the target style)
$ws.Range("D14").Value = 1
$ws.Range("F15").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").Value = -100
$ws.Range("H15").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null
$ws.Range("G14").Value = 1
$ws.Range("F15").Copy() | Out-Null
$ws.Range("G14").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Value = 100
$ws.Range("H15").Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null
$ws.Range("D23").Value = 1
$ws.Range("F15").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").Value = -100
$ws.Range("H15").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null

# --- Text (shared-string) updates that also require a style correction ---
# Set NumberFormat to text first so the literal string is not auto-coerced to a number,
# then restore the correct numeric-style formatting via PasteSpecial(xlPasteFormats).
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0